$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("13-11-2023", 63.85, 62.6, 66.4, 64.69, 192991570, 160, 28.5934, 7568, 10216, 357, 2547, 89),
  @("14-11-2023", 66.6, 62.35, 67.35, 65.35, 160602703, 160, 28.6341, 7671, 10656, 372, 2657, 93),
  @("15-11-2023", 67.25, 65.3, 68.5, 67.16, 133406071, 160, 28.6679, 7666, 10760, 375, 2682, 94),
  @("16-11-2023", 69.3, 65.45, 69.9, 68.45, 125135093, 160, 28.7045, 7781, 11088, 386, 2764, 96),
  @("17-11-2023", 68.45, 66.8, 71, 68.02, 106862508, 160, 28.7109, 7853, 10952, 381, 2730, 95),
  @("20-11-2023", 67, 67, 69.5, 68.05, 92416090, 160, 28.7568, 7946, 10720, 373, 2672, 93),
  @("21-11-2023", 68, 66.4, 69.1, 67.54, 101524367, 160, 28.7992, 8047, 10880, 378, 2712, 94),
  @("22-11-2023", 67.9, 67.35, 68.9, 68.01, 61613271, 160, 28.8136, 7987, 10864, 377, 2708, 94),
  @("23-11-2023", 70, 67.75, 71.5, 69.86, 171473631, 160, 28.8304, 7899, 11200, 388, 2792, 97),
  @("24-11-2023", 70.7, 68.75, 71.95, 70.94, 130470425, 160, 28.8607, 7960, 11312, 392, 2799, 97),
  @("27-11-2023", 70.85, 68.5, 71.65, 70.43, 158941993, 160, 28.9104, 8107, 11336, 392, 2805, 97),
  @("28-11-2023", 70.85, 69.7, 71.3, 70.4, 69847378, 160, 28.9289, 8094, 11336, 392, 2805, 97),
  @("29-11-2023", 70.1, 69.1, 70.8, 70.05, 52784958, 160, 28.9116, 8009, 11216, 388, 2775, 96),
  @("30-11-2023", 70.2, 67.7, 70.3, 68.99, 79839321, 160, 28.8728, 7949, 11232, 389, 2779, 96),
  @("01-12-2023", 70.35, 69.55, 72.5, 70.77, 74254032, 160, 28.9073, 8026, 11256, 389, 2785, 96),
  @("04-12-2023", 70.05, 68.4, 70.75, 69.45, 69668319, 160, 28.9035, 8087, 11208, 388, 2773, 96),
  @("05-12-2023", 70.5, 69.3, 71.2, 70.31, 55362651, 160, 28.9225, 8057, 11280, 390, 2791, 96),
  @("06-12-2023", 68.05, 67.2, 70.75, 69.52, 62510752, 160, 28.9313, 7855, 10888, 376, 2694, 93),
  @("07-12-2023", 69.75, 66, 69.75, 67.79, 63110923, 160, 28.926, 7979, 11160, 386, 2761, 95),
  @("08-12-2023", 71, 68.85, 71.8, 70.96, 73763656, 160, 28.9495, 7914, 11360, 392, 2810, 97),
  @("11-12-2023", 69.75, 66.15, 71.1, 68.7, 113052217, 160, 28.9917, 7728, 11160, 385, 2761, 95),
  @("12-12-2023", 70, 68.05, 70.2, 69.2, 45760584, 160, 29.0295, 7749, 11200, 386, 2771, 95),
  @("13-12-2023", 63, 63, 70, 65.49, 80748480, 160, 29.0579, 7529, 10080, 347, 2494, 86),
  @("14-12-2023", 63, 60.6, 65, 62.51, 86060442, 160, 29.0328, 7813, 10080, 347, 2494, 86),
  @("15-12-2023", 64.1, 62.6, 65.8, 64.39, 73162543, 160, 29.0369, 7994, 10256, 353, 2537, 87),
  @("18-12-2023", 65.5, 63.4, 65.75, 64.95, 58385390, 160, 29.0576, 7790, 10480, 361, 2593, 89),
  @("19-12-2023", 68, 62.95, 68.3, 66.98, 64345315, 160, 29.0822, 7726, 10880, 374, 2692, 93),
  @("20-12-2023", 70, 66.75, 71.25, 68.8, 69341393, 160, 29.1201, 7663, 11200, 385, 2771, 95),
  @("21-12-2023", 69.35, 67.45, 70.7, 69.33, 52906943, 160, 29.1422, 7765, 11096, 381, 2745, 94),
  @("22-12-2023", 72.5, 68.1, 73.95, 71.74, 195270455, 160, 29.1869, 7558, 11600, 397, 2870, 98)
)

$startRow = 48
$endRow = $startRow + $data.Count - 1

# Force column A (dates) to be treated as text so day<=12 strings
# (e.g. "01-12-2023") are not auto-parsed into date serial numbers.
$dateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$dateRange.NumberFormat = "@"

$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# Restore the default cell style on column A so the temporary text
# number format does not linger as a visible style difference.
$dateRange.Style = "Normal"
